$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: "freeday" moves from Tasks done (C5) to Additional info (D5)
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "freeday"

# New entries for chapter 7 / chapter 8 work
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "red chapter 7 completedfirst exercise, stuck on the second"

# Row 4: "chapter 4 was hard UwU" -> "chapter 4 was hard " (dropped the "UwU")
$ws.Range("D4").Value = "chapter 4 was hard "

$ws.Range("B11").Value = 0
$ws.Range("D11").Value = "internet died in the wholke area making it impossible to work"

$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "finished chapter 7"

$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "finished chapter 8"

# Update the active selection to match the final edit location
$ws.Range("C13").Select()
